$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Three new rows of benchmark data (re-run without back-to-front sorting).
# Row 20 reuses the existing "Unroll Plot" label; rows 21-22 introduce two
# brand-new labels ("oFast", "Combine uv mask").
$ws.Range("A20").Value2 = "Unroll Plot"
$ws.Range("B20").Value2 = 278

$ws.Range("A21").Value2 = "oFast"
$ws.Range("B21").Value2 = 279

$ws.Range("A22").Value2 = "Combine uv mask"
$ws.Range("B22").Value2 = 280

# Extend the C (FPS) and D (ratio vs baseline) formula columns down to the
# new rows, matching the existing B/30 and B/$B$2 pattern.
$ws.Range("C20").Formula = "=B20/30"
$ws.Range("C21").Formula = "=B21/30"
$ws.Range("C22").Formula = "=B22/30"

$ws.Range("D20").Formula = "=B20/`$B`$2"
$ws.Range("D21").Formula = "=B21/`$B`$2"
$ws.Range("D22").Formula = "=B22/`$B`$2"

$ws.Range("C20:C22").NumberFormat = $ws.Range("C19").NumberFormat
$ws.Range("D20:D22").NumberFormat = $ws.Range("D19").NumberFormat

$ws.Range("A22").Select()
